$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently has a header row (1) plus 4 identical data rows (2:5).
# Too many rows were being fetched to create a single req_order, so the
# fix duplicates that same data block to pad the sheet out to 10 data rows
# (rows 2:11), mirroring rows 2:5 into the newly inserted rows 6:11.
for ($i = 0; $i -lt 6; $i++) {
    $destRow = 6 + $i
    $srcRow = 2 + ($i % 4)
    $ws.Rows($srcRow).Copy()
    $ws.Rows($destRow).Insert()
}

# Move the active selection as recorded after the edit.
$ws.Range("C14").Select()
